$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B23").Value = "VUE.js program: Hogyan tudnám a CloseFun() funkcióban az isVisible értéket falsra állítani és így bezárni az ablakot?"
$ws.Range("B23").Select()
